$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the header row (row 1), shifting existing data down
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "julz0815"
$ws.Cells.Item(2, 2).Value = "dev-count-test-headless"
$ws.Cells.Item(2, 3).Value = "julz0815/dev-count-test-headless"

# "2025-12-05" looks like a date, so force the cell to text format before
# assigning it, otherwise Excel auto-converts it to a date serial number.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "2025-12-05"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).Value = "Y"
